$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.777.46"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "2.632.22"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'578.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("D6").Value = "'156.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("D7").Value = "'0.629"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -1.42%  "
$ws.Range("D10").Value = "'5.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("D13").Value = "'28.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.11%  "
$ws.Range("D14").Value = "3.107.54"
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("D16").Value = "63.703.74"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").Value = "2.651.84"
$ws.Range("E17").Value = "  +0.84%  "
$ws.Range("D18").Value = "'12.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("D19").Value = "'7.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.93%  "
$ws.Range("E20").Value = "  -2.43%  "
$ws.Range("D21").Value = "'344.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "'68.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.37%  "
$ws.Range("D24").Value = "'1.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.79%  "
$ws.Range("D25").Value = "'0.0000114"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.31%  "
$ws.Range("D26").Value = "'1.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.08%  "
$ws.Range("D27").Value = "'9.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("D28").Value = "'577.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.40%  "
$ws.Range("D29").Value = "'8.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.33%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("E33").Value = "  +1.39%  "
$ws.Range("D34").Value = "'6.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.80%  "
$ws.Range("D35").Value = "'5.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.28%  "
$ws.Range("E36").Value = "  -1.24%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "'0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "'19.71"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.90%  "
$ws.Range("D39").Value = "'1.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.82%  "
$ws.Range("D40").Value = "'153.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'2.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.57%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E43").Value = "  +3.74%  "
$ws.Range("D44").Value = "'24.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.90%  "
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("E46").Value = "  -1.45%  "
$ws.Range("D47").Value = "'0.633"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("E48").Value = "  -1.85%  "
$ws.Range("E49").Value = "  -0.96%  "
$ws.Range("E50").Value = "  +0.92%  "
$ws.Range("D51").Value = "'0.793"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.08%  "
